$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New row of sample-lineup data (row 15), following the same pattern/format
# as the existing rows (2-14). Copy the number formatting from row 14 into
# A15:E15 so the date/0.0/0.0% formats carry over correctly (column F is
# left with its existing column-level style, matching the other data rows).
foreach ($col in @("A", "B", "C", "D", "E")) {
    $ws.Range("$col`14").Copy() | Out-Null
    $ws.Range("$col`15").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
}
$excel.CutCopyMode = 0

$ws.Range("A15").Value = 43223
$ws.Range("B15").Value = 129.1
$ws.Range("C15").Formula = "=11219/38324"
$ws.Range("D15").Value = "Yes"
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 1
